$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_3_2_0"
$ws.Range("B2").Value = 0.4055803811774429
$ws.Range("C2").Value = 0.3283348312153617
$ws.Range("D2").Value = 0.2465232423111813
$ws.Range("E2").Value = 0.3364461700703654
$ws.Range("F2").Value = 0.6578474044799805
$ws.Range("G2").Value = 0.09549585729837418
$ws.Range("H2").Value = 1.003667593002319
$ws.Range("I2").Value = 0.5228708386421204

$ws.Range("A3").Value = "model_3_2_4"
$ws.Range("B3").Value = 0.4186187450236687
$ws.Range("C3").Value = -1.084735238539008
$ws.Range("D3").Value = 0.2411455610759831
$ws.Range("E3").Value = 0.1971876533440914
$ws.Range("F3").Value = 0.6434177756309509
$ws.Range("G3").Value = 0.2964030206203461
$ws.Range("H3").Value = 1.010830879211426
$ws.Range("I3").Value = 0.6326044797897339

$ws.Range("A4").Value = "model_3_2_1"
$ws.Range("B4").Value = 0.4190868080471527
$ws.Range("C4").Value = 0.3947613409603073
$ws.Range("D4").Value = 0.2243230112499317
$ws.Range("E4").Value = 0.3251308831232081
$ws.Range("F4").Value = 0.6428996920585632
$ws.Range("G4").Value = 0.08605148643255234
$ws.Range("H4").Value = 1.033239364624023
$ws.Range("I4").Value = 0.531787097454071

$ws.Range("A5").Value = "model_3_2_2"
$ws.Range("B5").Value = 0.4192054259238757
$ws.Range("C5").Value = 0.1860986654353012
$ws.Range("D5").Value = 0.2379263987058355
$ws.Range("E5").Value = 0.3160210564234603
$ws.Range("F5").Value = 0.6427684426307678
$ws.Range("G5").Value = 0.1157186850905418
$ws.Range("H5").Value = 1.015119075775146
$ws.Range("I5").Value = 0.5389655232429504

$ws.Range("A6").Value = "model_3_2_3"
$ws.Range("B6").Value = 0.4207452230514948
$ws.Range("C6").Value = -0.9126594501147605
$ws.Range("D6").Value = 0.240429410334262
$ws.Range("E6").Value = 0.213055052979268
$ws.Range("F6").Value = 0.6410642862319946
$ws.Range("G6").Value = 0.2719376683235168
$ws.Range("H6").Value = 1.011784791946411
$ws.Range("I6").Value = 0.6201012134552002

$ws.Range("A7").Value = "model_3_2_5"
$ws.Range("B7").Value = 0.4224606969271649
$ws.Range("C7").Value = -1.056744501057936
$ws.Range("D7").Value = 0.2457101950152081
$ws.Range("E7").Value = 0.2034924733220459
$ws.Range("F7").Value = 0.6391658782958984
$ws.Range("G7").Value = 0.2924233675003052
$ws.Range("H7").Value = 1.004750609397888
$ws.Range("I7").Value = 0.6276364326477051

$ws.Range("A8").Value = "model_3_2_6"
$ws.Range("B8").Value = 0.4323138953434896
$ws.Range("C8").Value = -1.23551690393924
$ws.Range("D8").Value = 0.2668015523441382
$ws.Range("E8").Value = 0.2031944373287931
$ws.Range("F8").Value = 0.6282612085342407
$ws.Range("G8").Value = 0.3178407847881317
$ws.Range("H8").Value = 0.9766558408737183
$ws.Range("I8").Value = 0.6278712153434753

$ws.Range("A9").Value = "model_3_2_7"
$ws.Range("B9").Value = 0.4389773484329662
$ws.Range("C9").Value = -1.35490193779806
$ws.Range("D9").Value = 0.2845573825768657
$ws.Range("E9").Value = 0.2059153281414184
$ws.Range("F9").Value = 0.6208868026733398
$ws.Range("G9").Value = 0.3348146975040436
$ws.Range("H9").Value = 0.9530042409896851
$ws.Range("I9").Value = 0.6257272362709045

$ws.Range("A10").Value = "model_3_2_8"
$ws.Range("B10").Value = 0.4409585839605559
$ws.Range("C10").Value = -1.360723298206814
$ws.Range("D10").Value = 0.2877810488903244
$ws.Range("E10").Value = 0.2079234659501225
$ws.Range("F10").Value = 0.6186941266059875
$ws.Range("G10").Value = 0.3356423676013947
$ws.Range("H10").Value = 0.9487101435661316
$ws.Range("I10").Value = 0.6241448521614075

$ws.Range("A11").Value = "model_3_2_9"
$ws.Range("B11").Value = 0.4498922820025838
$ws.Range("C11").Value = -1.557085848580881
$ws.Range("D11").Value = 0.3210803897013164
$ws.Range("E11").Value = 0.2156561261164551
$ws.Range("F11").Value = 0.6088070869445801
$ws.Range("G11").Value = 0.3635607957839966
$ws.Range("H11").Value = 0.9043538570404053
$ws.Range("I11").Value = 0.6180516481399536

$ws.Range("A12").Value = "model_3_2_10"
$ws.Range("B12").Value = 0.4551290195278714
$ws.Range("C12").Value = -1.612970047053316
$ws.Range("D12").Value = 0.3261117276168012
$ws.Range("E12").Value = 0.214320184817321
$ws.Range("F12").Value = 0.6030116677284241
$ws.Range("G12").Value = 0.3715062737464905
$ws.Range("H12").Value = 0.8976518511772156
$ws.Range("I12").Value = 0.6191043257713318

$ws.Range("A13").Value = "model_3_2_12"
$ws.Range("B13").Value = 0.4559440972838217
$ws.Range("C13").Value = -1.808603013565031
$ws.Range("D13").Value = 0.3317139932063585
$ws.Range("E13").Value = 0.2000897529241801
$ws.Range("F13").Value = 0.6021094918251038
$ws.Range("G13").Value = 0.3993209302425385
$ws.Range("H13").Value = 0.8901893496513367
$ws.Range("I13").Value = 0.6303176879882812

$ws.Range("A14").Value = "model_3_2_11"
$ws.Range("B14").Value = 0.4569522504727548
$ws.Range("C14").Value = -1.589079817770861
$ws.Range("D14").Value = 0.3255076392141771
$ws.Range("E14").Value = 0.2161219035859686
$ws.Range("F14").Value = 0.600993812084198
$ws.Range("G14").Value = 0.3681095838546753
$ws.Range("H14").Value = 0.8984565734863281
$ws.Range("I14").Value = 0.6176846027374268

$ws.Range("A15").Value = "model_3_2_18"
$ws.Range("B15").Value = 0.4589808447663438
$ws.Range("C15").Value = -1.798571465995723
$ws.Range("D15").Value = 0.3348859634723819
$ws.Range("E15").Value = 0.2035712378343991
$ws.Range("F15").Value = 0.5987487435340881
$ws.Range("G15").Value = 0.3978946506977081
$ws.Range("H15").Value = 0.8859641551971436
$ws.Range("I15").Value = 0.6275743842124939

$ws.Range("A16").Value = "model_3_2_13"
$ws.Range("B16").Value = 0.4592414429507937
$ws.Range("C16").Value = -1.671305234016091
$ws.Range("D16").Value = 0.330677875129532
$ws.Range("E16").Value = 0.2123803528772485
$ws.Range("F16").Value = 0.59846031665802
$ws.Range("G16").Value = 0.3798002302646637
$ws.Range("H16").Value = 0.8915694952011108
$ws.Range("I16").Value = 0.6206328272819519

$ws.Range("A17").Value = "model_3_2_14"
$ws.Range("B17").Value = 0.4597504637246015
$ws.Range("C17").Value = -1.678351614857075
$ws.Range("D17").Value = 0.3312593490185675
$ws.Range("E17").Value = 0.2121698971512195
$ws.Range("F17").Value = 0.5978970527648926
$ws.Range("G17").Value = 0.3808020353317261
$ws.Range("H17").Value = 0.8907949924468994
$ws.Range("I17").Value = 0.6207987070083618

$ws.Range("A18").Value = "model_3_2_15"
$ws.Range("B18").Value = 0.459875992586205
$ws.Range("C18").Value = -1.673481293093349
$ws.Range("D18").Value = 0.3314034313318789
$ws.Range("E18").Value = 0.2127497265380116
$ws.Range("F18").Value = 0.5977581739425659
$ws.Range("G18").Value = 0.3801096081733704
$ws.Range("H18").Value = 0.8906030654907227
$ws.Range("I18").Value = 0.6203417778015137

$ws.Range("A19").Value = "model_3_2_16"
$ws.Range("B19").Value = 0.4599591657116056
$ws.Range("C19").Value = -1.676359235807587
$ws.Range("D19").Value = 0.3316679486024275
$ws.Range("E19").Value = 0.2126853075317455
$ws.Range("F19").Value = 0.5976661443710327
$ws.Range("G19").Value = 0.3805187940597534
$ws.Range("H19").Value = 0.8902506828308105
$ws.Range("I19").Value = 0.6203925609588623

$ws.Range("A20").Value = "model_3_2_19"
$ws.Range("B20").Value = 0.4600576626711562
$ws.Range("C20").Value = -1.725154960326786
$ws.Range("D20").Value = 0.3310048048219237
$ws.Range("E20").Value = 0.2074964258261097
$ws.Range("F20").Value = 0.5975570678710938
$ws.Range("G20").Value = 0.387456476688385
$ws.Range("H20").Value = 0.8911340832710266
$ws.Range("I20").Value = 0.6244813203811646

$ws.Range("A21").Value = "model_3_2_17"
$ws.Range("B21").Value = 0.4602125046623332
$ws.Range("C21").Value = -1.680194748288858
$ws.Range("D21").Value = 0.3312115117538591
$ws.Range("E21").Value = 0.2119556128584603
$ws.Range("F21").Value = 0.5973857045173645
$ws.Range("G21").Value = 0.3810641169548035
$ws.Range("H21").Value = 0.8908587098121643
$ws.Range("I21").Value = 0.6209675669670105

$ws.Range("A22").Value = "model_3_2_20"
$ws.Range("B22").Value = 0.4607583715727566
$ws.Range("C22").Value = -1.717785460348948
$ws.Range("D22").Value = 0.3309621321192266
$ws.Range("E22").Value = 0.2081665582287074
$ws.Range("F22").Value = 0.5967816114425659
$ws.Range("G22").Value = 0.3864086866378784
$ws.Range("H22").Value = 0.8911908864974976
$ws.Range("I22").Value = 0.6239533424377441

$ws.Range("A23").Value = "model_3_2_21"
$ws.Range("B23").Value = 0.460830361179462
$ws.Range("C23").Value = -1.716744222749524
$ws.Range("D23").Value = 0.330269929615156
$ws.Range("E23").Value = 0.2077154725629082
$ws.Range("F23").Value = 0.596701979637146
$ws.Range("G23").Value = 0.3862606287002563
$ws.Range("H23").Value = 0.8921129107475281
$ws.Range("I23").Value = 0.6243087649345398

$ws.Range("A24").Value = "model_3_2_23"
$ws.Range("B24").Value = 0.4609464436864592
$ws.Range("C24").Value = -1.717186883308348
$ws.Range("D24").Value = 0.329047237081855
$ws.Range("E24").Value = 0.2067003480407607
$ws.Range("F24").Value = 0.5965734124183655
$ws.Range("G24").Value = 0.3863235712051392
$ws.Range("H24").Value = 0.8937416076660156
$ws.Range("I24").Value = 0.6251085996627808

$ws.Range("A25").Value = "model_3_2_24"
$ws.Range("B25").Value = 0.461240946151934
$ws.Range("C25").Value = -1.722601366952664
$ws.Range("D25").Value = 0.331034043143672
$ws.Range("E25").Value = 0.2077635624993396
$ws.Range("F25").Value = 0.5962475538253784
$ws.Range("G25").Value = 0.3870933949947357
$ws.Range("H25").Value = 0.8910951614379883
$ws.Range("I25").Value = 0.6242708563804626

$ws.Range("A26").Value = "model_3_2_22"
$ws.Range("B26").Value = 0.4613197048457837
$ws.Range("C26").Value = -1.704799133673089
$ws.Range("D26").Value = 0.3304327375078094
$ws.Range("E26").Value = 0.2089859670084752
$ws.Range("F26").Value = 0.596160352230072
$ws.Range("G26").Value = 0.3845622837543488
$ws.Range("H26").Value = 0.8918960094451904
$ws.Range("I26").Value = 0.6233076453208923
